$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.050.93"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "3.807.41"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "700.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("D7").Value = "3.807.98"
$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "4.453.05"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "3.820.52"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "71.177.52"
$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "512.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.25%  "

$ws.Range("E25").Value = "  -3.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.94%  "

$ws.Range("D27").Value = "3.959.40"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.02%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("E35").Value = "  -4.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("D37").Value = "3.777.73"
$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.11%  "

$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.71%  "

$ws.Range("E42").Value = "  -2.46%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.67%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "165.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "432.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000303"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.24%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.24%  "
